$p = $ppt.ActivePresentation
Write-Output ("HasTitleMaster before=" + $p.HasTitleMaster)
try {
    $tm = $p.AddTitleMaster()
    Write-Output ("AddTitleMaster ok: " + $tm)
} catch {
    Write-Output "AddTitleMaster failed: $_"
}
Write-Output ("HasTitleMaster after=" + $p.HasTitleMaster)
Write-Output ("Designs.Count=" + $p.Designs.Count)
